$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns C:H across rows 2..31 of the sheet.
# Row 2 gets brand-new sensor values; rows 3..21 inherit the previous
# occupant's C:H values (i.e. each row's C:H shifts down by one); rows
# 22..31 are new rows of data appended at the bottom. Columns A (timestamp)
# and B (label) are unaffected and already follow the existing pattern.
$chValues = @(
    @(-1.111309170722961, -16.8168830871582, -6.105263710021973, -0.4303629223038048, -0.01591238554787755, -0.3245061741155765),
    @(-0.0419791936874389, -20.49333000183105, -13.21517372131348, -1.073453715869358, 1.317456359841992, 1.647602679474009),
    @(-3.297537326812744, -12.87430763244629, -6.910844326019287, -0.4519367668809059, 0.3554307395169793, 2.567409327050211),
    @(-2.024067401885986, -10.72560882568359, -10.07154750823975, -0.7488680436831516, -0.1919719578713872, 1.722829146545486),
    @(2.292793273925781, -12.79408359527588, -4.549558639526367, 0.3624176147605125, 0.2589661164694446, 0.9208506445924814),
    @(11.05422019958496, -10.4705171585083, -6.454219818115234, -0.5349602806968903, 1.069952379880837, 1.569308005711615),
    @(14.52309226989746, -2.075542449951172, -4.196336269378662, -1.448478241928485, 0.2915999198160256, 2.148881192968673),
    @(65.91658020019531, -23.6422119140625, -20.2271671295166, 0.1628953875864145, 2.143472752150354, 0.7411342482356558),
    @(-11.09761142730713, -9.708492279052734, -16.66929626464844, 5.001680181056456, 0.5944304566423666, 0.3106197156325129),
    @(-24.99989891052246, -21.75775909423828, 0.6322221755981445, 2.196256959137806, 0.2997662000545932, -5.213745354854745),
    @(37.47088623046875, -10.36063957214356, 18.01900863647461, 1.67380741063286, -6.317273866729519, -8.133592685731504),
    @(-1.643178939819336, -5.723413467407227, -5.214177131652832, -3.399641977638681, -6.691358628393238, 2.842905701709222),
    @(10.17343330383301, -45.65782165527344, 20.15432167053223, -4.124252363413325, 2.333295746009847, 4.042209733434092),
    @(-1.11302375793457, -12.11679267883301, -19.5434455871582, -1.459144299270755, 10.65372753544007, 0.3427767938926811),
    @(25.58180236816406, -19.63376235961914, -23.0617904663086, 1.466152499703789, 0.4739527562085932, -0.06884271607679748),
    @(34.03413391113281, -47.92161560058594, -2.704976081848145, 7.103254157955869, -3.76272493550757, -1.125126261921481),
    @(-7.239111423492432, -3.941543340682984, -4.735682010650635, 6.002800260271336, 0.1307494660385649, -1.342715479746586),
    @(21.65508270263672, -27.37262153625488, 32.79093933105469, -0.9153044744699734, 0.02712116471860138, 1.211943073432978),
    @(-48.36841201782227, -30.63133239746094, -33.96612167358398, -9.937176976885045, 0.9889363561357609, 9.063289233616347),
    @(-3.451883316040039, -7.505331993103027, -13.28867816925049, 1.660898893820885, 0.03534786866495354, -1.336030611470736),
    @(-6.207223892211914, -16.8746509552002, -9.824567794799805, 3.122204536149511, -5.978971676165282, -1.097727024254715),
    @(-1.000881910324097, -4.219967842102051, 11.5488452911377, 5.814107340924815, -1.352176813518396, -1.502328406361931),
    @(26.53460311889648, 11.7469482421875, 5.673069000244141, -0.08294376405329906, -0.04769127909872317, 0.643018007278366),
    @(15.22074890136719, -62.06856155395508, 22.35076904296875, -4.552947707536836, -0.8177924519326552, 2.338970492867865),
    @(-22.60456275939941, 16.3265266418457, -2.586312532424927, -3.988243347456473, 3.882105143631221, -0.1545383128799389),
    @(-8.99936580657959, -18.44652366638184, -19.57757568359375, -1.513886836646995, 12.12395986789412, -6.341941402739815),
    @(80.43364715576172, -83.25542449951172, 37.27872467041016, 4.529875355608262, 0.5813608329838935, -4.973177166546146),
    @(-2.663472175598145, 1.494519710540772, -1.472611904144287, 5.40816806745133, -3.943533415554054, -4.557831178192283),
    @(18.3146915435791, -38.27093124389648, 47.83665084838867, 0.5592678154216362, 2.892510315951151, -1.382000362171883),
    @(20.0896053314209, -6.290483951568604, -35.30553817749023, -7.637950965837174, -7.75840919158016, -1.187092877235683)
)

$startRow = 2
for ($i = 0; $i -lt $chValues.Count; $i++) {
    $r = $startRow + $i
    $rowVals = $chValues[$i]
    for ($j = 0; $j -lt $rowVals.Count; $j++) {
        $col = 3 + $j
        $ws.Cells.Item($r, $col).Value = $rowVals[$j]
    }
}

# Fill in the timestamp (A) and label (B) columns for the newly appended rows.
for ($r = 22; $r -le 31; $r++) {
    $ws.Cells.Item($r, 1).Value = ($r - 2) * 100
    $ws.Cells.Item($r, 2).Value = "walkingToRunning"
}
